# "menys ifs a la GUI" - update VMIX sheet row 2 (SECTION 5 -> SECTION 2):
# shift the section-5 "final" block flags/countries/names by one slot, zero
# out the section points that no longer apply, swap the two leaders in the
# abbreviations row, and restyle the still-empty "-" placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VMIX")

# Section title
$ws.Range("F2").Value = "SECTION 2"

# Flags (F_BANDERA_1/2) and country codes (F_PAIS_1/2) and player names
# (F_PLAYER_1/2) swap position 1 <-> position 2.
$ws.Range("CO2").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\esp.png"
$ws.Range("CP2").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\fra.png"

$ws.Range("CU2").Value = "ESP"
$ws.Range("CV2").Value = "FRA"

$ws.Range("DA2").Value = "ALEJANDRO MO"
$ws.Range("DB2").Value = "VINCENT H"

# Section points for section 2 / 3 / 5 leaders reset to 0 (no longer
# applicable once the final standing recomputed).
$ws.Range("DG2").Value = 0
$ws.Range("DH2").Value = 0
$ws.Range("DM2").Value = 0
$ws.Range("DN2").Value = 0
$ws.Range("DT2").Value = 0
$ws.Range("DZ2").Value = 0
$ws.Range("EL2").Value = 0

# Abbreviation row: HER/MON swap.
$ws.Range("EQ2").Value = "MON"
$ws.Range("ER2").Value = "HER"

# 2_PUNTS_SECCIO / 2_PUNTS_P1 no longer scored.
$ws.Range("FD2").Value = 0
$ws.Range("FE2").Value = "-"

# The still-unplayed-section placeholders switch from "-" to " -" (leading
# space) across sections 3-6.
$dashRefs = @(
    "FL2","FM2","FN2","FO2","FP2","FQ2",
    "FS2","FT2","FU2","FV2","FW2","FX2",
    "FZ2","GA2","GB2","GC2","GD2","GE2",
    "GG2","GH2","GI2","GJ2","GK2","GL2"
)
foreach ($ref in $dashRefs) {
    $ws.Range($ref).Value = " -"
}
